$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.983.53'
$ws.Range("E2").Value = '  +0.92%  '

$ws.Range("D3").Value = '2.334.21'
$ws.Range("E3").Value = '  +4.59%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '97.70'
$ws.Range("E5").Value = '  +3.46%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '271.81'
$ws.Range("E6").Value = '  +0.90%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.628'
$ws.Range("E7").Value = '  +0.59%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.629'
$ws.Range("E9").Value = '  +0.93%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.99'
$ws.Range("E10").Value = '  -0.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0951'
$ws.Range("E11").Value = '  +2.88%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.17'
$ws.Range("E12").Value = '  -0.12%  '

$ws.Range("E13").Value = '  +0.42%  '

$ws.Range("D14").Value = '2.691.48'
$ws.Range("E14").Value = '  +4.80%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.67'

$ws.Range("E16").Value = '  +8.68%  '

$ws.Range("D17").Value = '2.340.68'
$ws.Range("E17").Value = '  +4.44%  '

$ws.Range("D18").Value = '43.914.72'
$ws.Range("E18").Value = '  +0.78%  '

$ws.Range("E19").Value = '  +5.95%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.44'
$ws.Range("E20").Value = '  +7.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.98'
$ws.Range("E21").Value = '  +3.53%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '240.17'
$ws.Range("E22").Value = '  +2.96%  '

$ws.Range("E23").Value = '  -2.00%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.48'
$ws.Range("E24").Value = '  +5.33%  '

$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.47'
$ws.Range("E26").Value = '  +2.25%  '

$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.53'
$ws.Range("E27").Value = '  +1.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.48'
$ws.Range("E28").Value = '  -1.79%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.28'
$ws.Range("E29").Value = '  +0.39%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '38.31'
$ws.Range("E30").Value = '  -5.24%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.53'
$ws.Range("E31").Value = '  +8.19%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.74'
$ws.Range("E32").Value = '  +0.40%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0906'
$ws.Range("E33").Value = '  -2.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.51'
$ws.Range("E34").Value = '  +0.35%  '

$ws.Range("E35").Value = '  +2.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0362'
$ws.Range("E36").Value = '  +3.64%  '

$ws.Range("E37").Value = '  -1.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.44'
$ws.Range("E38").Value = '  +2.62%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.41'
$ws.Range("E39").Value = '  -4.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.38'
$ws.Range("E40").Value = '  +8.92%  '

$ws.Range("E41").Value = '  +10.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.38'
$ws.Range("E42").Value = '  +19.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.24'
$ws.Range("E43").Value = '  -2.63%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.25'
$ws.Range("E44").Value = '  +10.50%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.74'
$ws.Range("E45").Value = '  -0.82%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.40'
$ws.Range("E46").Value = '  +1.23%  '

$ws.Range("E47").Value = '  +5.42%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '100.72'
$ws.Range("E48").Value = '  +0.11%  '

$ws.Range("E49").Value = '  +1.77%  '

$ws.Range("D50").Value = '2.568.53'
$ws.Range("E50").Value = '  +4.64%  '

$ws.Range("E51").Value = '  +15.00%  '

